# Repull data, push all data, mean calculation
# Updates the "dSF" column (F) values for several rows to reflect
# freshly repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -8
    3  = 7
    4  = -7
    5  = -2
    6  = 0
    7  = -1
    9  = 5
    10 = -2
    13 = 0
    18 = 1
    21 = 0
    22 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
